# Update the lottery-analysis workbook with a fresh run of results.
# Sheets (by name, in tab order):
#   1) Frequency Analysis   - Top 20 Numbers / Frequency Count
#   2) Suggested Numbers    - Suggested Numbers
#   3) Common Pairs         - Pair / Frequency
#   4) Consecutive Numbers  - Consecutive Sets
#   5) Range Analysis       - Range / Count
#   6) Hot Cold Analysis    - Hot Numbers / Cold Numbers

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Frequency Analysis: Top 20 Numbers with their frequency counts
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Frequency Analysis")
$freq = @(
    @(59,12), @(10,12), @(76,12), @(50,11), @(33,10), @(60,10), @(64,10), @(23,10),
    @(71,9),  @(18,8),  @(21,8),  @(28,8),  @(46,8),  @(26,8),  @(5,7),   @(31,7),
    @(47,7),  @(55,7),  @(56,7),  @(67,7)
)
for ($i = 0; $i -lt $freq.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $freq[$i][0]
    $ws1.Cells.Item($row, 2).Value = $freq[$i][1]
}

# ---------------------------------------------------------------------------
# 2) Suggested Numbers: same 20 numbers, in the same order
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Suggested Numbers")
$suggested = @(59, 10, 76, 50, 33, 60, 64, 23, 71, 18, 21, 28, 46, 26, 5, 31, 47, 55, 56, 67)
for ($i = 0; $i -lt $suggested.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $suggested[$i]
}

# ---------------------------------------------------------------------------
# 3) Common Pairs: pair label plus frequency (frequency stays 1 for each)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Common Pairs")
$pairs = @(
    "(50, 76)", "(64, 76)", "(33, 59)", "(7, 79)", "(10, 64)",
    "(10, 76)", "(28, 76)", "(60, 76)", "(10, 59)", "(21, 67)"
)
for ($i = 0; $i -lt $pairs.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 1).Value = $pairs[$i]
    $ws3.Cells.Item($row, 2).Value = 1
}

# ---------------------------------------------------------------------------
# 4) Consecutive Numbers: now 10 consecutive sets instead of 5
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Consecutive Numbers")
$consecutive = @(
    "((10, 11), 5)", "((55, 56), 4)", "((59, 60), 4)", "((25, 26), 4)", "((33, 34), 3)",
    "((73, 74), 3)", "((64, 65), 3)", "((22, 23), 3)", "((27, 28), 3)", "((9, 10), 3)"
)
for ($i = 0; $i -lt $consecutive.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 1).Value = $consecutive[$i]
}

# ---------------------------------------------------------------------------
# 5) Range Analysis: counts per quadrant (labels stay the same)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Range Analysis")
$ws5.Cells.Item(2, 2).Value = 117
$ws5.Cells.Item(3, 2).Value = 117
$ws5.Cells.Item(4, 2).Value = 119
$ws5.Cells.Item(5, 2).Value = 127

# ---------------------------------------------------------------------------
# 6) Hot Cold Analysis: hot numbers (from frequency top 10) paired with the
#    coldest numbers
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Hot Cold Analysis")
$hot = @(
    "(59, 12)", "(10, 12)", "(76, 12)", "(50, 11)", "(33, 10)",
    "(60, 10)", "(64, 10)", "(23, 10)", "(71, 9)", "(18, 8)"
)
$cold = @(
    "(52, 4)", "(54, 4)", "(20, 4)", "(58, 4)", "(78, 3)",
    "(36, 3)", "(45, 3)", "(38, 3)", "(30, 2)", "(43, 2)"
)
for ($i = 0; $i -lt $hot.Length; $i++) {
    $row = $i + 2
    $ws6.Cells.Item($row, 1).Value = $hot[$i]
    $ws6.Cells.Item($row, 2).Value = $cold[$i]
}
